$d = $word.ActiveDocument

$replacements = @(
    @{old="71×28="; new="17×57="},
    @{old="47×38="; new="70×11="},
    @{old="68×86="; new="86×91="},
    @{old="73×66="; new="13×68="},
    @{old="50×47="; new="79×28="},
    @{old="33×35="; new="92×49="},
    @{old="62×38="; new="34×33="},
    @{old="63×78="; new="74×48="},
    @{old="95×60="; new="32×17="},
    @{old="59×28="; new="50×22="},
    @{old="57×43="; new="13×46="},
    @{old="96×21="; new="91×35="},
    @{old="59×16="; new="89×12="},
    @{old="89×80="; new="37×89="},
    @{old="55×77="; new="12×54="},
    @{old="82×59="; new="94×32="},
    @{old="58×63="; new="62×56="},
    @{old="38×85="; new="83×89="},
    @{old="95×59="; new="91×68="},
    @{old="96×97="; new="73×14="},
    @{old="32×51="; new="40×79="},
    @{old="15×45="; new="86×11="},
    @{old="28×80="; new="69×67="},
    @{old="77×36="; new="89×89="},
    @{old="14×22="; new="28×98="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $r.new, 2)
}
